$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-detected as numbers by Excel (these cells store plain text in the
# workbook, e.g. "353.70", "41.53", etc).
$textForcedCells = @("D5", "D6", "D9", "D10", "D11", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D31", "D32", "D34", "D35", "D36", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update price (column D) and volume/change (column E) values row by row
$ws.Range("D2").Value = '51.820.89'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.822.91'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '353.70'
$ws.Range("E5").Value = '  +6.15%  '
$ws.Range("D6").Value = '112.76'
$ws.Range("E6").Value = '  -3.70%  '
$ws.Range("E7").Value = '  +4.71%  '
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  +4.55%  '
$ws.Range("D10").Value = '41.53'
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").Value = '0.0852'
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("E13").Value = '  -1.90%  '
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("D15").Value = '3.264.25'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '2.819.31'
$ws.Range("E16").Value = '  +1.88%  '
$ws.Range("D17").Value = '0.886'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '51.717.93'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").Value = '7.43'
$ws.Range("E19").Value = '  +8.52%  '
$ws.Range("D20").Value = '3.20'
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").Value = '13.41'
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").Value = '0.0₃0991'
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("D23").Value = '269.71'
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").Value = '69.74'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +3.82%  '
$ws.Range("D26").Value = '26.75'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("D31").Value = '50.70'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").Value = '33.97'
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("E33").Value = '  +26.17%  '
$ws.Range("D34").Value = '5.85'
$ws.Range("E34").Value = '  +4.49%  '
$ws.Range("D35").Value = '5.34'
$ws.Range("E35").Value = '  +6.57%  '
$ws.Range("D36").Value = '0.0824'
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  -1.74%  '
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").Value = '18.19'
$ws.Range("E40").Value = '  -5.14%  '
$ws.Range("D41").Value = '23.94'
$ws.Range("E41").Value = '  +3.01%  '
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").Value = '  +1.95%  '
$ws.Range("D43").Value = '126.18'
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").Value = '2.51'
$ws.Range("E44").Value = '  +1.76%  '
$ws.Range("D45").Value = '2.30'
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").Value = '2.090.48'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  +0.40%  '
$ws.Range("D48").Value = '2.27'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").Value = '5.66'
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("E50").Value = '  +6.50%  '
$ws.Range("D51").Value = '60.61'
$ws.Range("E51").Value = '  -0.33%  '

# Restore the default (Normal) style on the cells we temporarily reformatted
# as text, so no stray formatting is left behind on those cells.
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
